# Covid_Anx.xlsx "Add files via upload" edit
#
# The cells in "feeling safe" (B2:E3) and "work affected" (B2:C2) contained a
# literal placeholder "x" (shared string). That placeholder text is cleared
# out, leaving the cells blank (formatting/style is kept). Because that "x"
# shared-string becomes unused, it disappears from the shared-string table on
# save and every subsequent shared string shifts down by one index - that
# part is automatic bookkeeping handled by the engine, not something this
# script needs to touch directly.
#
# The other, purely cosmetic, parts of the change: the lingering multi-cell
# selection (A1:A8) left over from when the sheets were first authored gets
# reset - to A1 on most sheets, and to B2 (the first cell that was edited) on
# the two sheets that were actually touched. The "marital status" sheet also
# gets explicit (wider) column widths for columns C:E instead of the
# worksheet-wide default.

$wb = $excel.ActiveWorkbook

$wsLoneliness   = $wb.Worksheets.Item("loneliness")
$wsSex          = $wb.Worksheets.Item("sex")
$wsMarital      = $wb.Worksheets.Item("marital status")
$wsFeelingSafe  = $wb.Worksheets.Item("feeling safe")
$wsWorkAffected = $wb.Worksheets.Item("work affected")
$wsDisability   = $wb.Worksheets.Item("disability")

# "feeling safe": clear the placeholder "x" values, keep the cell styling.
$wsFeelingSafe.Range("B2:E3").ClearContents()
$wsFeelingSafe.Rows.Item(3).AutoFit()

# "work affected": same placeholder cleanup.
$wsWorkAffected.Range("B2:C2").ClearContents()

# "marital status": give columns C, D, E explicit widths instead of the
# sheet-wide default.
$wsMarital.Columns.Item(3).ColumnWidth = 34.45
$wsMarital.Columns.Item(4).ColumnWidth = 28.06
$wsMarital.Columns.Item(5).ColumnWidth = 22.09

# Reset/update the lingering A1:A8 selections left on each sheet.
$wsLoneliness.Range("A1").Select()
$wsSex.Range("A1").Select()
$wsMarital.Range("A1").Select()
$wsFeelingSafe.Range("B2").Select()
$wsWorkAffected.Range("B2").Select()
$wsDisability.Range("A1").Select()
